$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New fixture list (match, date | competition) replacing the previous
# Norwich/Sporting-away block with a new Fulham FA Cup tie, shifting
# later fixtures up, updating the Man United game date and adding a
# new Sporting CP home Champions League tie.
$fixtures = @(
    @("Manchester City v Fulham", "05/02/2022 15:00 | The Emirates FA Cup"),
    @("Manchester City v Brentford", "09/02/2022 19:45 | Premier League"),
    @("Manchester City v Tottenham Hotspur", "19/02/2022 17:30 | Premier League"),
    @("Manchester City v Manchester United", "06/03/2022 16:30 | Premier League"),
    @("Manchester City v Sporting CP", "09/03/2022 20:00 | UEFA Champions League"),
    @("Manchester City v Brighton & Hove Albion ", "19/03/2022 15:00 | Premier League"),
    @("Manchester City v Liverpool ", "09/04/2022 15:00 | Premier League"),
    @("Manchester City v Watford", "23/04/2022 15:00 | Premier League"),
    @("Manchester City v Newcastle United", "07/05/2022 15:00 | Premier League"),
    @("Manchester City v Aston Villa", "22/05/2022 15:00 | Premier League")
)

for ($i = 0; $i -lt $fixtures.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $fixtures[$i][0]
    $ws.Cells.Item($row, 2).Value = $fixtures[$i][1]
}
